# Update the "through" date in the sheet name and the header label,
# then refresh the carjacking counts that changed with the new day of data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Sheet tab name: Through 2022-06-08 -> Through 2022-06-09
$ws.Name = "Through 2022-06-09"

# Header cell B1 label: "June 2022 (through June 08)" -> "June 2022 (through June 09)"
$ws.Range("B1").Value = "June 2022 (through June 09)"

# New / updated carjacking counts by neighborhood-month
$ws.Range("T3").Value = 1
$ws.Range("AF5").Value = 3
$ws.Range("AL6").Value = 1
$ws.Range("T9").Value = 1
$ws.Range("Z12").Value = 1
$ws.Range("T19").Value = 1
$ws.Range("B20").Value = 2
$ws.Range("AF30").Value = 1
$ws.Range("AF31").Value = 1
$ws.Range("AF40").Value = 1
$ws.Range("B47").Value = 1
$ws.Range("AF55").Value = 1
$ws.Range("B57").Value = 1
$ws.Range("B66").Value = 2
$ws.Range("B89").Value = 1
$ws.Range("N94").Value = 3
